$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark that currently sits at the end of
#    the "WriteState2" paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Replace the text of the last list item in the "A faire :" section.
$d.Content.Find.Execute(
    "Faire la première passe en ne changeant que les paramètres nécessaires (lineaire, deform, Dphi ?)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Tout mettre dans une fonction", 2)

# 3. Append a new list paragraph "Ajouter une interpolation" after it, keeping
#    the same list formatting (InsertParagraphAfter carries over pPr).
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()
$interpPara = $d.Paragraphs.Last
$interpPara.Range.Text = "Ajouter une interpolation"

# 4. Append one more (empty) list paragraph that will host the relocated
#    "_GoBack" bookmark, again preserving the list formatting.
$interpPara.Range.InsertParagraphAfter()
$finalPara = $d.Paragraphs.Last

# Temporarily give the new paragraph a single placeholder character so we can
# anchor the bookmark to a proper (non-collapsed) range, then remove the
# placeholder text again -- this leaves the paragraph with just the bookmark
# and no run, matching Word's own output for an empty _GoBack target.
$finalPara.Range.Text = "Z"
$placeholderRange = $d.Range($finalPara.Range.Start, $finalPara.Range.Start + 1)
$d.Bookmarks.Add("_GoBack", $placeholderRange)
$deleteRange = $d.Range($finalPara.Range.Start, $finalPara.Range.Start + 1)
$deleteRange.Delete()
